$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @(3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,0.3952656666666667,1.185797,0.02469993039301214,0.02469993039301213,7.887835739595112,70.990521656356,0.0002899115295515849,0.0002899115295515848)
$data[3] = @(3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,3.484068333333333,10.452205,0.2177174811148058,0.2177174811148058,69.52731045581555,625.7457941023399,0.002555424527753674,0.002555424527753673)
$data[4] = @(3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,0.5492906666666667,1.647872,0.03432486647933305,0.03432486647933305,10.96152516482845,98.653726483456,0.0004028826957946675,0.0004028826957946675)
$data[5] = @(3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,11.574079,34.722237,0.7232577220128491,0.723257722012849,230.9698050908307,2078.728245817476,0.008489123212592571,0.00848912321259257)
$data[6] = @(3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,0.3952656666666667,1.185797,0.02469993039301214,0.02469993039301213,647.1856081648789,5824.67047348391,0.0237868251521784,0.0237868251521784)
$data[7] = @(3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,3.484068333333333,10.452205,0.2177174811148058,0.2177174811148058,5704.616093301794,51341.54483971615,0.2096689170150749,0.2096689170150749)
$data[8] = @(3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,0.5492906666666667,1.647872,0.03432486647933305,0.03432486647933305,899.3774166217955,8094.39674959616,0.03305594729719381,0.03305594729719381)
$data[9] = @(3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,11.574079,34.722237,0.7232577220128491,0.723257722012849,18950.74120586412,170556.6708527771,0.6965203828408231,0.696520382840823)
$data[10] = @(3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,0.3952656666666667,1.185797,0.02469993039301214,0.02469993039301213,6.91747302100289,62.25725718902601,0.0002542465703340967,0.0002542465703340968)
$data[11] = @(3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,3.484068333333333,10.452205,0.2177174811148058,0.2177174811148058,60.97405044665445,548.76645401989,0.002241055824630099,0.002241055824630099)
$data[12] = @(3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,0.5492906666666667,1.647872,0.03432486647933305,0.03432486647933305,9.613036718819558,86.517330469376,0.0003533200070075979,0.000353320007007598)
$data[13] = @(3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,11.574079,34.722237,0.7232577220128491,0.723257722012849,202.5558655287273,1823.002789758546,0.007444790020195426,0.007444790020195426)
$data[14] = @(3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,0.3952656666666667,1.185797,0.02469993039301214,0.02469993039301213,10.03821562009889,90.34394058088999,0.0003689471409480521,0.000368947140948052)
$data[15] = @(3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,3.484068333333333,10.452205,0.2177174811148058,0.2177174811148058,88.48182909509444,796.3364618558498,0.003252083747347087,0.003252083747347087)
$data[16] = @(3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,0.5492906666666667,1.647872,0.03432486647933305,0.03432486647933305,13.94985351651555,125.54868164864,0.0005127164793369762,0.0005127164793369762)
$data[17] = @(3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,11.574079,34.722237,0.7232577220128491,0.723257722012849,293.9367377537433,2645.43063978369,0.01080342593923805,0.01080342593923805)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}
